# Add "tno"/"tvol"/"tval" (trade count/volume/value) and "z" (share count)
# variable rows to the RealTimeVariables table.
#
# The table in column A/B is a flat list of variable-name -> Persian
# description pairs. We insert:
#   - 3 new rows right before the existing "bvol" row: tno, tvol, tval
#   - 1 new row right after "bvol" (before "buy_i_count"): z
# All subsequent rows shift down accordingly; no existing data changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Locate the "bvol" row (variable name in column A) so this keeps working
# even if the sheet layout shifts slightly.
$bvolRow = $ws.Range("A1:A1000").Find("bvol").Row

# Insert 3 rows above "bvol" for tno / tvol / tval.
# Fill column A first (top to bottom), then column B (top to bottom) -
# matches how the authoring session entered the new variable names before
# filling in their Persian descriptions.
$ws.Rows($bvolRow + ":" + ($bvolRow + 2)).Insert()
$ws.Range("A" + $bvolRow).Value = "tno"
$ws.Range("A" + ($bvolRow + 1)).Value = "tvol"
$ws.Range("A" + ($bvolRow + 2)).Value = "tval"
$ws.Range("B" + $bvolRow).Value = "تعداد معاملات"
$ws.Range("B" + ($bvolRow + 1)).Value = "حجم معاملات"
$ws.Range("B" + ($bvolRow + 2)).Value = "ارزش معاملات"

# "bvol" itself is now 3 rows further down.
$bvolRow = $bvolRow + 3

# Insert 1 row right after "bvol" for z.
$zRow = $bvolRow + 1
$ws.Rows($zRow + ":" + $zRow).Insert()
$ws.Range("A" + $zRow).Value = "z"
$ws.Range("B" + $zRow).Value = "تعداد سهام"

# Leave the selection where the edits ended up.
$ws.Range("C1048576").Select()
